$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change the base "Damage" min-stat input value from 2 to 10
$ws.Range("L39").Value = 10

# Update the active view: top-left visible cell and the current selection
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("L40").Select()
